$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial data for rows 2-6 (2014/12 - 2018/12 columns) with corrected values
# Row 2
$ws.Range("D2").Value = 2506
$ws.Range("E2").Value = -116
$ws.Range("F2").Value = -116
$ws.Range("G2").Value = -299
$ws.Range("H2").Value = -243
$ws.Range("I2").Value = -244
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5046
$ws.Range("L2").Value = 3427
$ws.Range("M2").Value = 1619
$ws.Range("N2").Value = 1599
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 642
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = -34
$ws.Range("S2").Value = 8
$ws.Range("T2").Value = 70
$ws.Range("U2").Value = -65
$ws.Range("V2").Value = 2509
$ws.Range("W2").Value = -4.61
$ws.Range("X2").Value = -9.7
$ws.Range("Y2").Value = -18.07
$ws.Range("Z2").Value = -5.08
$ws.Range("AA2").Value = 211.58
$ws.Range("AB2").Value = 57.74
$ws.Range("AC2").Value = -1582
$ws.Range("AD2").Value = -3.92
$ws.Range("AE2").Value = 9746
$ws.Range("AF2").Value = 0.64
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 16578736

# Row 3
$ws.Range("D3").Value = 2448
$ws.Range("E3").Value = -464
$ws.Range("F3").Value = -464
$ws.Range("G3").Value = -746
$ws.Range("H3").Value = -593
$ws.Range("I3").Value = -593
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 5547
$ws.Range("L3").Value = 4371
$ws.Range("M3").Value = 1176
$ws.Range("N3").Value = 1157
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 753
$ws.Range("Q3").Value = -97
$ws.Range("R3").Value = -185
$ws.Range("S3").Value = 578
$ws.Range("T3").Value = 122
$ws.Range("U3").Value = -218
$ws.Range("V3").Value = 3023
$ws.Range("W3").Value = -18.96
$ws.Range("X3").Value = -24.22
$ws.Range("Y3").Value = -43.06
$ws.Range("Z3").Value = -11.19
$ws.Range("AA3").Value = 371.61
$ws.Range("AB3").Value = -26.77
$ws.Range("AC3").Value = -3328
$ws.Range("AD3").Value = -1.71
$ws.Range("AE3").Value = 6006
$ws.Range("AF3").Value = 0.95
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 19443468

# Row 4
$ws.Range("D4").Value = 4219
$ws.Range("E4").Value = -223
$ws.Range("F4").Value = -223
$ws.Range("G4").Value = -448
$ws.Range("H4").Value = -371
$ws.Range("I4").Value = -372
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 5419
$ws.Range("L4").Value = 4232
$ws.Range("M4").Value = 1187
$ws.Range("N4").Value = 1166
$ws.Range("O4").Value = 21
$ws.Range("P4").Value = 1122
$ws.Range("Q4").Value = -322
$ws.Range("R4").Value = -69
$ws.Range("S4").Value = 211
$ws.Range("T4").Value = 56
$ws.Range("U4").Value = -377
$ws.Range("V4").Value = 2878
$ws.Range("W4").Value = -5.29
$ws.Range("X4").Value = -8.79
$ws.Range("Y4").Value = -32.03
$ws.Range("Z4").Value = -6.77
$ws.Range("AA4").Value = 356.38
$ws.Range("AB4").Value = -44.84
$ws.Range("AC4").Value = -1721
$ws.Range("AD4").Value = -2.66
$ws.Range("AE4").Value = 4380
$ws.Range("AF4").Value = 1.05
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 26803263

# Row 5
$ws.Range("D5").Value = 5565
$ws.Range("E5").Value = -394
$ws.Range("F5").Value = -394
$ws.Range("G5").Value = -1105
$ws.Range("H5").Value = -1073
$ws.Range("I5").Value = -1062
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 6615
$ws.Range("L5").Value = 6050
$ws.Range("M5").Value = 565
$ws.Range("N5").Value = 557
$ws.Range("O5").Value = 15
$ws.Range("P5").Value = 271
$ws.Range("Q5").Value = -27
$ws.Range("R5").Value = -1388
$ws.Range("S5").Value = 1502
$ws.Range("T5").Value = 660
$ws.Range("U5").Value = -687
$ws.Range("V5").Value = 3947
$ws.Range("W5").Value = -7.08
$ws.Range("X5").Value = -19.29
$ws.Range("Y5").Value = -123.19
$ws.Range("Z5").Value = -17.84
$ws.Range("AA5").Value = 1070.7
$ws.Range("AB5").Value = -100.35
$ws.Range("AC5").Value = -3531
$ws.Range("AD5").Value = -0.35
$ws.Range("AE5").Value = 1031
$ws.Range("AF5").Value = 1.19
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 54216657

# Row 6
$ws.Range("D6").Value = 4123
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 22
$ws.Range("G6").Value = -271
$ws.Range("H6").Value = -381
$ws.Range("I6").Value = -381
$ws.Range("K6").Value = 5884
$ws.Range("L6").Value = 5606
$ws.Range("M6").Value = 279
$ws.Range("N6").Value = 278
$ws.Range("P6").Value = 367
$ws.Range("Q6").Value = -198
$ws.Range("R6").Value = -158
$ws.Range("S6").Value = 395
$ws.Range("T6").Value = 387
$ws.Range("U6").Value = -586
$ws.Range("V6").Value = 4216
$ws.Range("W6").Value = 0.53
$ws.Range("X6").Value = -9.23
$ws.Range("Y6").Value = -91.06
$ws.Range("Z6").Value = -6.09
$ws.Range("AA6").Value = 2012.65
$ws.Range("AB6").Value = -172.39
$ws.Range("AC6").Value = -699
$ws.Range("AD6").Value = -1.34
$ws.Range("AE6").Value = 380
$ws.Range("AF6").Value = 2.46
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 73374054

# Remove rows 7-9 numeric data (2019/12(E) - 2021/12(E)) - columns D through AI,
# keeping only the row label columns (A, B, C)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
